# New crime data collected — weekly CompStat refresh (60th Precinct).
# Moves the reporting week forward one week (Volume/Number + date range)
# and refreshes every crime-count / percent-change figure in the
# "Crime Complaints" table (rows 14-30, cols C:N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats — used below to copy a donor cell's number format/style
# onto a cell whose content is switching between text ("n/a"-style values
# stored as shared strings "0" / "***.*") and numeric, without disturbing
# the value we just wrote into it.
$xlPasteFormats = -4122

# ---- Masthead: volume/number and the reporting week text ----
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# ---- Row 14: Murder ----
# F14 flips from a numeric 1 to the literal "0" placeholder text used
# elsewhere in the sheet for not-applicable counts.
$ws.Range("F14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial($xlPasteFormats)

# ---- Row 15: Rape ----
$ws.Range("L15").Value = -50
$ws.Range("M15").Value = 100
$ws.Range("L14").Copy()
$ws.Range("M15").PasteSpecial($xlPasteFormats)

# ---- Row 16: Robbery ----
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 8.333333333333
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = 14
$ws.Range("L16").Value = 128
$ws.Range("M16").Value = -37.362637362637
$ws.Range("N16").Value = -83.852691218130

# ---- Row 17: Fel. Assault ----
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 12.5
$ws.Range("F17").Value = 23
$ws.Range("H17").Value = -4.166666666666
$ws.Range("I17").Value = 117
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = 15.841584158415
$ws.Range("L17").Value = 1.739130434782
$ws.Range("M17").Value = 134
$ws.Range("N17").Value = -53.386454183266

# ---- Row 18: Burglary ----
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = -21.951219512195
$ws.Range("L18").Value = -30.434782608695
$ws.Range("M18").Value = -33.333333333333
$ws.Range("N18").Value = -91.061452513966

# ---- Row 19: Gr. Larceny ----
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -11.111111111111
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 11.764705882352
$ws.Range("I19").Value = 174
$ws.Range("J19").Value = 149
$ws.Range("K19").Value = 16.778523489932
$ws.Range("L19").Value = 65.714285714285
$ws.Range("M19").Value = 16.778523489932
$ws.Range("N19").Value = -7.446808510638

# ---- Row 20: G.L.A. ----
# D20/E20 flip from numeric to the "0" / "***.*" placeholder text.
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial($xlPasteFormats)
$ws.Range("E20").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial($xlPasteFormats)
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 14.285714285714
$ws.Range("I20").Value = 31
$ws.Range("K20").Value = -22.5
$ws.Range("L20").Value = 93.75
$ws.Range("M20").Value = -16.216216216216
$ws.Range("N20").Value = -91.644204851752

# ---- Row 21: TOTAL ----
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 28.571428571428
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = 1.190476190476
$ws.Range("I21").Value = 416
$ws.Range("J21").Value = 385
$ws.Range("K21").Value = 8.051948051948
$ws.Range("L21").Value = 32.484076433121
$ws.Range("M21").Value = 10.344827586206
$ws.Range("N21").Value = -72.934287573194

# ---- Row 22: Transit ----
# D22/E22 flip from the "0" / "***.*" placeholder text to real numbers.
$ws.Range("D22").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial($xlPasteFormats)
$ws.Range("E22").Value = -100
$ws.Range("L14").Copy()
$ws.Range("E22").PasteSpecial($xlPasteFormats)
$ws.Range("G22").Value = 3
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -66.666666666666
$ws.Range("M22").Value = -82.352941176470

# ---- Row 23: Housing ----
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 17
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 41.666666666666
$ws.Range("I23").Value = 51
$ws.Range("J23").Value = 49
$ws.Range("K23").Value = 4.081632653061
$ws.Range("L23").Value = -25
$ws.Range("M23").Value = 82.142857142857

# ---- Row 24: Petit Larceny ----
$ws.Range("C24").Value = 27
$ws.Range("E24").Value = 22.727272727272
$ws.Range("F24").Value = 105
$ws.Range("G24").Value = 95
$ws.Range("H24").Value = 10.526315789473
$ws.Range("I24").Value = 367
$ws.Range("J24").Value = 407
$ws.Range("K24").Value = -9.828009828009
$ws.Range("L24").Value = 39.543726235741
$ws.Range("M24").Value = 10.210210210210

# ---- Row 25: Misd. Assault ----
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 62.857142857142
$ws.Range("I25").Value = 189
$ws.Range("J25").Value = 167
$ws.Range("K25").Value = 13.173652694610
$ws.Range("L25").Value = 36.956521739130
$ws.Range("M25").Value = 27.702702702702

# ---- Row 26: UCR Rape* ----
# C26 flips from numeric 1 to the "0" placeholder text.
$ws.Range("C26").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C26").PasteSpecial($xlPasteFormats)
$ws.Range("L26").Value = -37.5

# ---- Row 27: Other Sex Crimes ----
# D27/E27 flip from numeric to the "0" / "***.*" placeholder text.
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$ws.Range("E27").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial($xlPasteFormats)
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 21
$ws.Range("K27").Value = 31.25
$ws.Range("L27").Value = 75

# ---- Row 30: Hate Crimes ----
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = 25
$ws.Range("L30").Value = 400
